$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: last-updated timestamp string
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Agosto de 2020 a las 17:27"

# Updated COVID-19 statistics per country (columns: B=Casos totales, C=Nuevos casos,
# D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)
$ws.Cells.Item(4, 2).Value = 5754059
$ws.Cells.Item(4, 3).Value = 7787
$ws.Cells.Item(4, 4).Value = 3096527
$ws.Cells.Item(4, 5).Value = 2479939
$ws.Cells.Item(4, 7).Value = 169
$ws.Cells.Item(4, 8).Value = 177593
$ws.Cells.Item(5, 2).Value = 3505361
$ws.Cells.Item(5, 3).Value = 264
$ws.Cells.Item(5, 5).Value = 739509
$ws.Cells.Item(5, 7).Value = 22
$ws.Cells.Item(5, 8).Value = 112445
$ws.Cells.Item(15, 2).Value = 323313
$ws.Cells.Item(15, 3).Value = 1033
$ws.Cells.Item(15, 7).Value = 2
$ws.Cells.Item(15, 8).Value = 41405
$ws.Cells.Item(16, 4).Value = 239806
$ws.Cells.Item(16, 5).Value = 74511
$ws.Cells.Item(16, 7).Value = 50
$ws.Cells.Item(16, 8).Value = 6567
$ws.Cells.Item(20, 2).Value = 257065
$ws.Cells.Item(20, 3).Value = 947
$ws.Cells.Item(20, 4).Value = 204960
$ws.Cells.Item(20, 5).Value = 16678
$ws.Cells.Item(20, 7).Value = 9
$ws.Cells.Item(20, 8).Value = 35427
$ws.Cells.Item(22, 4).Value = 205800
$ws.Cells.Item(22, 5).Value = 16264
$ws.Cells.Item(33, 4).Value = 77785
$ws.Cells.Item(33, 5).Value = 21005
$ws.Cells.Item(33, 7).Value = 14
$ws.Cells.Item(33, 8).Value = 809
$ws.Cells.Item(59, 5).Value = 3432
$ws.Cells.Item(59, 7).Value = 2
$ws.Cells.Item(59, 8).Value = 2000
$ws.Cells.Item(65, 2).Value = 32484
$ws.Cells.Item(65, 3).Value = 547
$ws.Cells.Item(65, 4).Value = 22683
$ws.Cells.Item(65, 5).Value = 8872
$ws.Cells.Item(65, 7).Value = 8
$ws.Cells.Item(65, 8).Value = 929
$ws.Cells.Item(66, 2).Value = 31763
$ws.Cells.Item(66, 3).Value = 322
$ws.Cells.Item(66, 4).Value = 18157
$ws.Cells.Item(66, 5).Value = 13074
$ws.Cells.Item(66, 7).Value = 16
$ws.Cells.Item(66, 8).Value = 532
$ws.Cells.Item(147, 2).Value = 1406
$ws.Cells.Item(147, 3).Value = 11
$ws.Cells.Item(147, 5).Value = 508
$ws.Cells.Item(158, 1).Value = "Lesoto"
$ws.Cells.Item(158, 2).Value = 1015
$ws.Cells.Item(158, 3).Value = 19
$ws.Cells.Item(158, 4).Value = 472
$ws.Cells.Item(158, 5).Value = 513
$ws.Cells.Item(158, 8).Value = 30
$ws.Cells.Item(159, 1).Value = "Vietnam"
$ws.Cells.Item(159, 2).Value = 1009
$ws.Cells.Item(159, 3).Value = 2
$ws.Cells.Item(159, 4).Value = 545
$ws.Cells.Item(159, 5).Value = 439
$ws.Cells.Item(159, 8).Value = 25
$ws.Cells.Item(164, 2).Value = 806
$ws.Cells.Item(164, 3).Value = 39
$ws.Cells.Item(164, 4).Value = 162
$ws.Cells.Item(164, 5).Value = 632
$ws.Cells.Item(172, 1).Value = "Birmania"
$ws.Cells.Item(172, 2).Value = 419
$ws.Cells.Item(172, 3).Value = 20
$ws.Cells.Item(172, 4).Value = 337
$ws.Cells.Item(172, 5).Value = 76
$ws.Cells.Item(172, 8).Value = 6
$ws.Cells.Item(173, 1).Value = "Comoras"
$ws.Cells.Item(173, 2).Value = 417
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 396
$ws.Cells.Item(173, 5).Value = 14
$ws.Cells.Item(173, 8).Value = 7
$ws.Cells.Item(174, 2).Value = 384
$ws.Cells.Item(174, 3).Value = 1
$ws.Cells.Item(174, 4).Value = 300
$ws.Cells.Item(174, 5).Value = 84
$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(213, 4).Value = 12
$ws.Cells.Item(213, 8).Value = 1
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0

Write-Host "Update complete"

